# Update cryptocurrency price/volume data (and two row label swaps)
# per the scraped GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.144.48'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.44%  '
$ws.Range('E2').NumberFormat = 'General'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.657.90'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('E3').NumberFormat = 'General'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('E4').NumberFormat = 'General'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.58'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5161'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.97%  '
$ws.Range('E6').NumberFormat = 'General'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2635'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.22%  '
$ws.Range('E8').NumberFormat = 'General'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06274'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.12%  '
$ws.Range('E9').NumberFormat = 'General'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.74'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -5.03%  '
$ws.Range('E10').NumberFormat = 'General'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07723'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.94%  '
$ws.Range('E11').NumberFormat = 'General'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('B12').NumberFormat = 'General'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('C12').NumberFormat = 'General'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.430'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.75%  '
$ws.Range('E12').NumberFormat = 'General'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('B13').NumberFormat = 'General'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C13').NumberFormat = 'General'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.653.40'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.96%  '
$ws.Range('E13').NumberFormat = 'General'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.884.57'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.36%  '
$ws.Range('E14').NumberFormat = 'General'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5411'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.21%  '
$ws.Range('E15').NumberFormat = 'General'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8115'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.74%  '
$ws.Range('E16').NumberFormat = 'General'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.54%  '
$ws.Range('E17').NumberFormat = 'General'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.166.32'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('E18').NumberFormat = 'General'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.613'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.78%  '
$ws.Range('E20').NumberFormat = 'General'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.43'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.79%  '
$ws.Range('E21').NumberFormat = 'General'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.07'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.53%  '
$ws.Range('E22').NumberFormat = 'General'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.999'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -5.16%  '
$ws.Range('E23').NumberFormat = 'General'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '139.73'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.58%  '
$ws.Range('E25').NumberFormat = 'General'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1222'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.90%  '
$ws.Range('E26').NumberFormat = 'General'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.185'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.19%  '
$ws.Range('E28').NumberFormat = 'General'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.427'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.78%  '
$ws.Range('E29').NumberFormat = 'General'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05967'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -5.59%  '
$ws.Range('E30').NumberFormat = 'General'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.40%  '
$ws.Range('E31').NumberFormat = 'General'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.557'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.37%  '
$ws.Range('E32').NumberFormat = 'General'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.254'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.39%  '
$ws.Range('E33').NumberFormat = 'General'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.34%  '
$ws.Range('E34').NumberFormat = 'General'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9628'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.11%  '
$ws.Range('E35').NumberFormat = 'General'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.427'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('E36').NumberFormat = 'General'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.772'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5677'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -7.71%  '
$ws.Range('E38').NumberFormat = 'General'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01588'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.66%  '
$ws.Range('E39').NumberFormat = 'General'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.953'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.49%  '
$ws.Range('E40').NumberFormat = 'General'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.04%  '
$ws.Range('E41').NumberFormat = 'General'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.22%  '
$ws.Range('E42').NumberFormat = 'General'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.007.40'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -8.12%  '
$ws.Range('E43').NumberFormat = 'General'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.54'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('E44').NumberFormat = 'General'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.799.16'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.37%  '
$ws.Range('E45').NumberFormat = 'General'
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'Aave'
$ws.Range('B46').NumberFormat = 'General'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('C46').NumberFormat = 'General'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.71'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.28%  '
$ws.Range('E46').NumberFormat = 'General'
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('B47').NumberFormat = 'General'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('C47').NumberFormat = 'General'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₈105'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.84%  '
$ws.Range('E47').NumberFormat = 'General'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.64%  '
$ws.Range('E48').NumberFormat = 'General'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.988'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.40%  '
$ws.Range('E49').NumberFormat = 'General'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05170'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.448'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.42%  '
$ws.Range('E51').NumberFormat = 'General'

Write-Output "Updated $(91) cells"
